# Actualización automática 2025-08-27 17:00:09
#
# A new salesperson row ("MORAN MARQUEZ DAYSE MARCELA") is inserted right
# above "MOROCHO PLAZA SHIRLEY AURELIA" on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, pushing all the following rows (and the trailing
# summary row) down by one. The new row carries zeros for every metric
# column, and the summary row's "X de 20" style labels become "X de 21"
# on "VENTAS POR GRUPO" to reflect the new headcount.

$wb = $excel.ActiveWorkbook

$newName = "MORAN MARQUEZ DAYSE MARCELA"
$office  = "OFICINA-CATAECSA"
$anchor  = "MOROCHO PLAZA SHIRLEY AURELIA"

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": columns A..R, metrics in C..R
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Find the row whose CLIENTE (col B) currently holds the anchor name.
$insertRow1 = 0
for ($r = 2; $r -le 30; $r++) {
    if ($ws1.Cells.Item($r, 2).Value() -eq $anchor) {
        $insertRow1 = $r
        break
    }
}

$ws1.Rows.Item($insertRow1).Insert()
$ws1.Cells.Item($insertRow1, 1).Value = $office
$ws1.Cells.Item($insertRow1, 2).Value = $newName
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item($insertRow1, $col).Value = 0
}

# Last row is now the "X de N" summary row - bump the denominator by one.
$lastRow1 = $ws1.UsedRange.Rows.Count
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item($lastRow1, $col)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = $v.Replace(" de 20", " de 21")
    }
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": columns A..G, metrics in C..G
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$insertRow2 = 0
for ($r = 2; $r -le 30; $r++) {
    if ($ws2.Cells.Item($r, 2).Value() -eq $anchor) {
        $insertRow2 = $r
        break
    }
}

$ws2.Rows.Item($insertRow2).Insert()
$ws2.Cells.Item($insertRow2, 1).Value = $office
$ws2.Cells.Item($insertRow2, 2).Value = $newName
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item($insertRow2, $col).Value = 0
}

Write-Output "done"
